$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B to fit new content (target stored width 16.6640625;
# COM ColumnWidth snaps to the nearest pixel-width grid step, so use the
# character-width value that rounds to the closest match)
$ws.Columns.Item(2).ColumnWidth = 15.83

# Add new row of data (row 6)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "jenkins/jenkins"
$ws.Range("C6").Value = "6a44d1dd2d60"
$ws.Range("D6").Value = "jenikins"
$ws.Range("E6").Value = "cfb509230b4e"
$ws.Range("F6").Value = "lts"

# Update the selected cell to match the post-edit state
$ws.Range("A7").Select()
